$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$carrera = "Tecnicatura Universitaria en Procesamiento y Explotación de Datos"

$ws.Range("A84").Value = $carrera
$ws.Range("B84").Value = "Algoritmos y Estructuras de Datos - TUPED"
$ws.Range("C84").Value = "Aho, A. V., Hopcroft, J. E., & Ullman, J. D. (1998). Estructuras de datos y algoritmos (A. Vargas Villazón & J. Lozano Moreno, Trads.; 1.a ed.). Pearson Educación."

$ws.Range("A85").Value = $carrera
$ws.Range("B85").Value = "Algoritmos y Estructuras de Datos - TUPED"
$ws.Range("C85").Value = "Bhargava, A. (2016). Grokking Algorithms: An Illustrated Guide for Programmers and Other Curious People."

$ws.Range("A86").Value = $carrera
$ws.Range("B86").Value = "Algoritmos y Estructuras de Datos - TUPED"
$ws.Range("C86").Value = "Cormen, T. H., Leiserson, C. E., Rivest, R. L., & Stein, C. (2022). Introduction to algorithms (Fourth edition). The MIT Press."

$ws.Range("A87").Value = $carrera
$ws.Range("B87").Value = "Algoritmos y Estructuras de Datos - TUPED"
$ws.Range("C87").Value = "Kok, A. S. (2019). Hands-On Blockchain for Python Developers: Gain blockchain programming skills to build decentralized applications using Python. Packt Publishing Ltd."

$ws.Range("A88").Value = $carrera
$ws.Range("B88").Value = "Algoritmos y Estructuras de Datos - TUPED"
$ws.Range("C88").Value = "Miller, B., & Ranum, D. (2013). Solución de problemas con algoritmos y estructuras de datos usando Python (M. Orozco-Alzate, Trad.; 2.a ed.)."

$ws.Range("A89").Value = $carrera
$ws.Range("B89").Value = "Aspectos Legales del Uso de la Información - TUPED"
$ws.Range("C89").Value = "“Tratado de Derecho Constitucional”, ROSATTI Horacio, 2° edición ampliada, Rubinzal Culzoni Editores, Santa Fé 2017."

$ws.Range("A90").Value = $carrera
$ws.Range("B90").Value = "Aspectos Legales del Uso de la Información - TUPED"
$ws.Range("C90").Value = "“Código Civil y Comercial de la Nación comentado”, dirigido por Ricardo Luis LORENZETTI, 1° edición, Rubinzal Culzoni Editores, Santa Fé, 2014."

$ws.Range("A91").Value = $carrera
$ws.Range("B91").Value = "Aspectos Legales del Uso de la Información - TUPED"
$ws.Range("C91").Value = "“Derecho Procesal Administrativo”, HUTCHINSON Tomás, 1° edición, Rubinzal Culzoni Editores, Santa Fé 2009."

$ws.Range("A92").Value = $carrera
$ws.Range("B92").Value = "Aspectos Legales del Uso de la Información - TUPED"
$ws.Range("C92").Value = "“Acerca de la reflexión humana. La necesidad de detenernos y partir de la ignorancia”, Guillermo MAGI, Fundación La Hendija, 1° edición, 2016."

$ws.Range("A93").Value = $carrera
$ws.Range("B93").Value = "Aspectos Legales del Uso de la Información - TUPED"
$ws.Range("C93").Value = "“Tratado de derecho administrativo y obras selectas”, GORDILLO Agustín, Tomo 1, Parte general, 11a ed., Buenos Aires, F.D.A., 2013."

$ws.Range("A94").Value = $carrera
$ws.Range("B94").Value = "Aspectos Legales del Uso de la Información - TUPED"
$ws.Range("C94").Value = "Desregulación, Entre el Derecho y la Economía, Jorge Eduardo BUSTAMANTE, Editorial Abeledo Perrot, Buenos Aires, 1993."

$ws.Range("A95").Value = $carrera
$ws.Range("B95").Value = "Aspectos Legales del Uso de la Información - TUPED"
$ws.Range("C95").Value = "Sistema económico y rentístico, ALBERDI, Juan B., Editorial Ciudad Argentina, Buenos Aires, 1998."

$ws.Range("A96").Value = $carrera
$ws.Range("B96").Value = "Aspectos Legales del Uso de la Información - TUPED"
$ws.Range("C96").Value = "“Manual De Derecho Procesal Civil”, Lino E. Palacio Actualizador: Carlos E. Camps , Luis E. Palacio , Lino A. Palacio , Editorial Abeledo Perrot, Edicion 2016."

$ws.Range("A97").Value = $carrera
$ws.Range("B97").Value = "Aspectos Legales del Uso de la Información - TUPED"
$ws.Range("C97").Value = "“Contratos administrativos”, SILVA CENSIO, Jorge A., Astrea, Buenos Aires, 1982."

$ws.Range("A98").Value = $carrera
$ws.Range("B98").Value = "Aspectos Legales del Uso de la Información - TUPED"
$ws.Range("C98").Value = "Nuevos Principios de Comercio Internacional, Para actuar en Escenarios Globalizados, LEDESMA Carlos A., Ediciones Macchi, 5° edición, 1997."

$ws.Range("A99").Value = $carrera
$ws.Range("B99").Value = "Probabilidad y Estadística - TUPED"
$ws.Range("C99").Value = "* MENDENHALL, W.: Introducción a la probabilidad y estadística., México DF: Cengage Learning Editores, SA de CV (2010)."

$ws.Range("A100").Value = $carrera
$ws.Range("B100").Value = "Probabilidad y Estadística - TUPED"
$ws.Range("C100").Value = "* RAMOS, EVA: Estadística para todos. Ediciones Pirámide (2016)"

$ws.Range("A101").Value = $carrera
$ws.Range("B101").Value = "Probabilidad y Estadística - TUPED"
$ws.Range("C101").Value = "* GUTIERREZ BANEGAS, A. L. y SABARIA, L.: Probabilidad y estadística: enfoque por competencias. McGraw - Hill (2012)"

$ws.Range("A102").Value = $carrera
$ws.Range("B102").Value = "Probabilidad y Estadística - TUPED"
$ws.Range("C102").Value = "McGraw - Hill (2012)"

$ws.Range("C84").WrapText = $true

# Update the view so the active cell / visible area reflect where the author
# was working when the new rows were added.
$ws.Range("C86").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 1
